# Redis.NetCore command-tracking sheet update.
# Marks several Server-group commands as "Finished" (column C) and records
# their implementing method name (column E) or a comment (column F),
# matching commit: "INFO, BGSAVE, SAVE, CONFIG GET/SET/RESETSTATS/REWRITE,
# BGREWRITEAOF, LASTSAVE, DBSIZE, FLUSHDB, FLUSHALL".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# row -> (Finished?, column letter for the note, note text)
$updates = @(
    @{ Row = 82;  Finished = $true;  Col = "E"; Text = "BackgroundRewriteAppendOnlyFileAsync" },  # BRREWRITEAOF
    @{ Row = 83;  Finished = $true;  Col = "E"; Text = "BackgroundSaveAsync" },                     # BGSAVE
    @{ Row = 94;  Finished = $true;  Col = "E"; Text = "GetConfigurationAsync" },                   # CONFIG GET
    @{ Row = 95;  Finished = $true;  Col = "E"; Text = "RewriteConfigurationAsync" },                # CONFIG REWRITE
    @{ Row = 96;  Finished = $true;  Col = "E"; Text = "SetConfigurationAsync" },                   # CONFIG SET
    @{ Row = 97;  Finished = $true;  Col = "E"; Text = "ResetConfigurationStatsAsync" },            # CONFIG RESETSTAT
    @{ Row = 98;  Finished = $true;  Col = "E"; Text = "GetDatabaseSizeAsync" },                    # DBSIZE
    @{ Row = 101; Finished = $true;  Col = "E"; Text = "FlushAllAsync" },                           # FLUSHALL
    @{ Row = 102; Finished = $true;  Col = "E"; Text = "FlushDatabaseAsync" },                      # FLUSHDB
    @{ Row = 103; Finished = $true;  Col = "E"; Text = "GetServerInformationAsync" },               # INFO
    @{ Row = 104; Finished = $true;  Col = "E"; Text = "GetLastSaveDateTimeAsync" },                # LASTSAVE
    @{ Row = 105; Finished = $false; Col = "F"; Text = "No plans, requires streaming" },            # MONITOR
    @{ Row = 107; Finished = $true;  Col = "E"; Text = "SaveAsync" },                                # SAVE
    @{ Row = 108; Finished = $false; Col = "F"; Text = "No plans" }                                  # SHUTDOWN
)

foreach ($u in $updates) {
    if ($u.Finished) {
        $ws.Cells.Item($u.Row, 3).Value = $true
    }
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Text
}

# Restore the current selection/view state recorded in the sheet.
$ws.Range("C120").Select()
